$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking values so they stay as text (matches source inlineStr cells)
$forceTextCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D17", "D19", "D20", "D23", "D24", "D26", "D29", "D30", "D32", "D33", "D34", "D39", "D40", "D49", "D50", "D51")
foreach ($c in $forceTextCells) { $ws.Range($c).NumberFormat = "@" }

# Row 2
$ws.Range("D2").Value = '49.998.20'
$ws.Range("E2").Value = '  +3.87%  '

# Row 3
$ws.Range("D3").Value = '2.659.56'
$ws.Range("E3").Value = '  +6.19%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").Value = '328.43'
$ws.Range("E5").Value = '  +2.26%  '

# Row 6
$ws.Range("D6").Value = '111.09'
$ws.Range("E6").Value = '  +3.18%  '

# Row 7
$ws.Range("D7").Value = '0.530'
$ws.Range("E7").Value = '  +0.88%  '

# Row 8
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("D9").Value = '0.561'
$ws.Range("E9").Value = '  +3.99%  '

# Row 10
$ws.Range("D10").Value = '40.70'
$ws.Range("E10").Value = '  +2.69%  '

# Row 11
$ws.Range("D11").Value = '20.68'
$ws.Range("E11").Value = '  +2.60%  '

# Row 12
$ws.Range("D12").Value = '0.0823'
$ws.Range("E12").Value = '  +1.27%  '

# Row 13
$ws.Range("E13").Value = '  +0.80%  '

# Row 14
$ws.Range("E14").Value = '  +2.89%  '

# Row 15
$ws.Range("D15").Value = '3.073.87'
$ws.Range("E15").Value = '  +6.17%  '

# Row 16
$ws.Range("D16").Value = '2.629.86'
$ws.Range("E16").Value = '  +5.26%  '

# Row 17
$ws.Range("D17").Value = '0.883'
$ws.Range("E17").Value = '  +5.78%  '

# Row 18
$ws.Range("D18").Value = '49.931.94'
$ws.Range("E18").Value = '  +4.01%  '

# Row 19
$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D19").Value = '13.32'
$ws.Range("E19").Value = '  +2.61%  '

# Row 20
$ws.Range("B20").Value = 'ImmutableX'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D20").Value = '3.02'
$ws.Range("E20").Value = '  +8.44%  '

# Row 21
$ws.Range("E21").Value = '  +1.89%  '

# Row 22
$ws.Range("E22").Value = '  +2.52%  '

# Row 23
$ws.Range("D23").Value = '282.50'
$ws.Range("E23").Value = '  +1.80%  '

# Row 24
$ws.Range("D24").Value = '73.31'
$ws.Range("E24").Value = '  +2.51%  '

# Row 25
$ws.Range("E25").Value = '  +2.34%  '

# Row 26
$ws.Range("D26").Value = '27.01'
$ws.Range("E26").Value = '  +3.86%  '

# Row 27
$ws.Range("E27").Value = '  +0.02%  '

# Row 28
$ws.Range("E28").Value = '  +6.84%  '

# Row 29
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = '10.00'
$ws.Range("E29").Value = '  +2.80%  '

# Row 30
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '36.73'
$ws.Range("E30").Value = '  +3.79%  '

# Row 31
$ws.Range("E31").Value = '  +2.49%  '

# Row 32
$ws.Range("D32").Value = '49.74'
$ws.Range("E32").Value = '  +0.15%  '

# Row 33
$ws.Range("D33").Value = '19.59'
$ws.Range("E33").Value = '  +0.03%  '

# Row 34
$ws.Range("D34").Value = '5.45'
$ws.Range("E34").Value = '  +2.64%  '

# Row 35
$ws.Range("E35").Value = '  -0.14%  '

# Row 37
$ws.Range("E37").Value = '  +6.76%  '

# Row 38
$ws.Range("E38").Value = '  +2.65%  '

# Row 39
$ws.Range("D39").Value = '3.12'
$ws.Range("E39").Value = '  +8.29%  '

# Row 40
$ws.Range("D40").Value = '125.62'
$ws.Range("E40").Value = '  +3.52%  '

# Row 41
$ws.Range("E41").Value = '  +1.69%  '

# Row 42
$ws.Range("E42").Value = '  +5.15%  '

# Row 43
$ws.Range("E43").Value = '  +0.90%  '

# Row 44
$ws.Range("E44").Value = '  +3.69%  '

# Row 45
$ws.Range("E45").Value = '  +7.08%  '

# Row 46
$ws.Range("D46").Value = '2.070.11'
$ws.Range("E46").Value = '  +2.21%  '

# Row 47
$ws.Range("E47").Value = '  +14.10%  '

# Row 48
$ws.Range("E48").Value = '  +8.52%  '

# Row 49
$ws.Range("D49").Value = '9.10'
$ws.Range("E49").Value = '  +1.16%  '

# Row 50
$ws.Range("D50").Value = '5.42'
$ws.Range("E50").Value = '  +4.61%  '

# Row 51
$ws.Range("D51").Value = '81.81'

# Remove the temporary style markers so cells end up with no style attribute (matches source)
foreach ($c in $forceTextCells) { $ws.Range($c).ClearFormats() }